# Biology Practicals - gen3 Oct29 Light: add Male/Female sub-totals
# (summed across the vial rows) to each of the four sheets, matching
# the "Data analysis for Biology completed" commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1 ("Vial 2 G3 / (Vial 1 gone!)")
#   Totals for Red/White across rows 3,6,9 (Male) and 4,7,10 (Female)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("G2").Value = "Male"
$ws1.Range("H2").Value = "Female"

$ws1.Range("F3").Value = "Red"
$ws1.Range("G3").Formula = "=SUM(B3,B6,B9)"
$ws1.Range("H3").Formula = "=SUM(C3,C6,C9)"

$ws1.Range("F4").Value = "White"
$ws1.Range("G4").Formula = "=SUM(B4,B7,B10)"
$ws1.Range("H4").Formula = "=SUM(C4,C7,C10)"

# ---------------------------------------------------------------------
# Sheet2 ("Vial 3")
#   Totals for Red/White across rows 3,6 (Male) and 4,7 (Female)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("G2").Value = "Male"
$ws2.Range("H2").Value = "Female"

$ws2.Range("F3").Value = "Red"
$ws2.Range("G3").Formula = "=SUM(B3,B6)"
$ws2.Range("H3").Formula = "=SUM(C3,C6)"

$ws2.Range("F4").Value = "White"
$ws2.Range("G4").Formula = "=SUM(B4,B7)"
$ws2.Range("H4").Formula = "=SUM(C4,C7)"

# ---------------------------------------------------------------------
# Sheet3 ("Vial 4") - summary lands in columns H:J (I/J hold the totals)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("I2").Value = "Male"
$ws3.Range("J2").Value = "Female"

$ws3.Range("H3").Value = "Red"
$ws3.Range("I3").Formula = "=SUM(B3,B6)"
$ws3.Range("J3").Formula = "=SUM(C3,C6)"

$ws3.Range("H4").Value = "White"
$ws3.Range("I4").Formula = "=SUM(B4,B7)"
$ws3.Range("J4").Formula = "=SUM(C4,C7)"

# ---------------------------------------------------------------------
# Sheet4 ("Vial 5")
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("G2").Value = "Male"
$ws4.Range("H2").Value = "Female"

$ws4.Range("F3").Value = "Red"
$ws4.Range("G3").Formula = "=SUM(B3,B6)"
$ws4.Range("H3").Formula = "=SUM(C3,C6)"

$ws4.Range("F4").Value = "White"
$ws4.Range("G4").Formula = "=SUM(B4,B7)"
$ws4.Range("H4").Formula = "=SUM(C4,C7)"

# ---------------------------------------------------------------------
# View state: selections per sheet, and Sheet4 becomes the active tab
# (matching the author having finished work there).
# ---------------------------------------------------------------------
$ws1.Range("F2:H2").Select()
$ws2.Range("F5").Select()
$ws3.Range("J2").Select()
$ws4.Range("H2").Select()

$ws4.Activate()
